# fakeleaderboard.xlsx update:
#  - replace the placeholder crypto leaderboard with a "who scored what"
#    leaderboard (rank / name / score)
#  - give the top score (C2) a thousands-separator number format
#  - move the cell selection to F8
#  - nudge the saved window position (cosmetic; mirrors the author's
#    workbookView xWindow/yWindow change)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "rank"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "score"

# Row 2 - Alex, 12000 (formatted with a thousands separator)
$ws.Range("B2").Value = "Alex"
$ws.Range("C2").Value = 12000
$ws.Range("C2").NumberFormat = "#,##0"

# Row 3 - Shimron, 11000
$ws.Range("B3").Value = "Shimron"
$ws.Range("C3").Value = 11000

# Row 4 - Roi, 10000
$ws.Range("B4").Value = "Roi"
$ws.Range("C4").Value = 10000

# Row 5 - Daniel, 9000
$ws.Range("B5").Value = "Daniel"
$ws.Range("C5").Value = 9000

# Note: row 7 (Bill Gates) is entered before row 6 (Elon Mask) so the
# shared-string table ends up in the same order the workbook was authored in.
# Row 7 - Bill Gates, 7000
$ws.Range("B7").Value = "Bill Gates"
$ws.Range("C7").Value = 7000

# Row 6 - Elon Mask, 8000
$ws.Range("B6").Value = "Elon Mask"
$ws.Range("C6").Value = 8000

# Row 8 - Mark Zuckerberg, 6000
$ws.Range("B8").Value = "Mark Zuckerberg"
$ws.Range("C8").Value = 6000

# Row 9 - Jeff Bezos, 5000
$ws.Range("B9").Value = "Jeff Bezos"
$ws.Range("C9").Value = 5000

# Row 10 - Steve Jobs, 0
$ws.Range("B10").Value = "Steve Jobs"
$ws.Range("C10").Value = 0

# Row 11 - Drake, -2000
$ws.Range("B11").Value = "Drake"
$ws.Range("C11").Value = -2000

# Move the selection, as in the diff's sheetView/selection change
$ws.Range("F8").Select() | Out-Null

# Move the workbook window (matches the xWindow/yWindow change on
# bookViews/workbookView)
$win = $wb.Windows.Item(1)
$win.Left = 31530
$win.Top = 7035
